# Remove the click-to-reveal animation timing (p:timing) from slide 7.
# The slide's speaker notes cover self-attention for tax document
# processing; the author dropped the per-bullet "appear" animation
# sequence (11 click effects targeting shape id 6) that PowerPoint
# stores as the slide's <p:timing> element.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

$seq = $s.TimeLine.MainSequence
for ($i = $seq.Count; $i -ge 1; $i--) {
    $seq.Item($i).Delete()
}
